$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

# --- Per-row price (D) / 1h volume change (E) updates ---
Set-TextValue "D2" "26.939.83"
Set-TextValue "E2" "  -0.53%  "
Set-TextValue "D3" "1.665.79"
Set-TextValue "E3" "  +0.57%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "215.71"
Set-TextValue "E5" "  +0.29%  "
Set-TextValue "E6" "  +4.56%  "
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "E8" "  +0.93%  "
Set-TextValue "E9" "  +0.07%  "
Set-TextValue "E10" "  +2.96%  "
Set-TextValue "D11" "0.0896"
Set-TextValue "E11" "  +3.74%  "
Set-TextValue "E12" "  +0.54%  "
Set-TextValue "D13" "1.665.57"
Set-TextValue "E13" "  +0.44%  "
Set-TextValue "D14" "4.08"
Set-TextValue "E14" "  -0.05%  "
Set-TextValue "E15" "  +0.62%  "
Set-TextValue "D16" "66.16"
Set-TextValue "E16" "  +1.69%  "
Set-TextValue "D17" "26.912.68"
Set-TextValue "D18" "234.65"
Set-TextValue "E18" "  -1.58%  "
Set-TextValue "D19" "8.02"
Set-TextValue "E19" "  +0.46%  "
Set-TextValue "E20" "  +0.35%  "
Set-TextValue "E21" "  +0.10%  "
Set-TextValue "E22" "  -2.04%  "
Set-TextValue "D23" "2.22"
Set-TextValue "E23" "  -1.07%  "
Set-TextValue "E24" "  -1.55%  "
Set-TextValue "D25" "146.26"
Set-TextValue "E25" "  +0.29%  "
Set-TextValue "D26" "7.13"
Set-TextValue "E26" "  -0.42%  "
Set-TextValue "D27" "0.116"
Set-TextValue "E27" "  +1.56%  "
Set-TextValue "D28" "15.89"
Set-TextValue "E28" "  +0.29%  "
Set-TextValue "D29" "0.999"
Set-TextValue "E29" "  -0.03%  "
Set-TextValue "D30" "0.0497"
Set-TextValue "E30" "  -0.07%  "
Set-TextValue "E31" "  +0.10%  "
Set-TextValue "E32" "  +2.18%  "
Set-TextValue "D33" "1.460.53"
Set-TextValue "E33" "  -4.51%  "
Set-TextValue "D34" "3.15"
Set-TextValue "E34" "  +2.56%  "
Set-TextValue "D35" "1.64"
Set-TextValue "E35" "  +2.84%  "
Set-TextValue "E36" "  -0.25%  "
Set-TextValue "D37" "0.583"
Set-TextValue "E37" "  +0.65%  "
Set-TextValue "D38" "0.906"
Set-TextValue "E38" "  +1.73%  "
Set-TextValue "E39" "  -0.04%  "
Set-TextValue "D40" "5.72"
Set-TextValue "E40" "  -3.79%  "
Set-TextValue "E41" "  +0.04%  "
Set-TextValue "E42" "  +0.61%  "
Set-TextValue "D43" "0.978"
Set-TextValue "E43" "  +6.54%  "
Set-TextValue "E44" "  -1.05%  "
Set-TextValue "D45" "1.809.15"
Set-TextValue "E45" "  +0.67%  "
Set-TextValue "D46" "0.782"
Set-TextValue "E46" "  +0.73%  "
Set-TextValue "D47" "90.39"
Set-TextValue "E47" "  +0.24%  "
Set-TextValue "E50" "  +3.89%  "
Set-TextValue "E51" "  +0.17%  "

# --- Rows 48/49 swapped: BabyDogeCoin now ranks above RenderToken ---
Set-TextValue "B48" "BabyDogeCoin"
Set-TextValue "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.0₆0106"
Set-TextValue "E48" "  +0.14%  "
Set-TextValue "B49" "RenderToken"
Set-TextValue "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.54"
Set-TextValue "E49" "  +0.58%  "
